# "Generate Report for Handback"
#
# Marks the 6cf4087f... handback as complete:
#   - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#     on the Overview sheet (zh-cn/de-de summary columns) and on the per-language
#     "Status" column of the zh-cn / de-de sheets.
#   - The "Latest Target File"(I)/"Latest Handback File"(J)/"Latest Handback
#     DateTime"(K) columns get populated for the 6cf4087f... row (row 2) and the
#     fffff5a89d48... row (row 3) on both the zh-cn and de-de sheets.
#   - A couple of columns get widened so the new, longer values are readable.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$hoMdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3200212d1817e5c43515c21ef2fac59cd81d9583/e2e/6cf4087f-bfdd-4afa-b937-c2e3a4e86a9a.md"
$hoMdDisplay = "6cf4087f-bfdd-4afa-b937-c2e3a4e86a9a.md"

# Widths below are expressed in the `Range.ColumnWidth` character-unit scale;
# the host stores/rounds widths to the nearest 1/6th of a character (i.e. a
# whole screen pixel with the default font), so these are chosen to land on
# the desired rendered widths (~30 and 40 characters respectively).
$wideStatusColWidth = 29.166666666666668
$wideFileColWidth   = 39.166666666666664

# ---------------------------------------------------------------------------
# Overview sheet: widen the zh-cn / de-de status columns and refresh their
# status text for both rows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E1:F1").ColumnWidth = $wideStatusColWidth

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C1").ColumnWidth = $wideStatusColWidth
$wsZhCn.Range("I1:J1").ColumnWidth = $wideFileColWidth

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $hoMdUrl, "", "", $hoMdDisplay)
$wsZhCn.Range("J2").Value = "6cf4087f-bfdd-4afa-b937-c2e3a4e86a9a.3dd1741d221c505698bb1d476ff36879f498027c.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-25 07:02:01"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $hoMdUrl, "", "", $hoMdDisplay)
$wsZhCn.Range("J3").Value = "6cf4087f-bfdd-4afa-b937-c2e3a4e86a9a.3dd1741d221c505698bb1d476ff36879f498027c.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-25 07:02:01"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C1").ColumnWidth = $wideStatusColWidth
$wsDeDe.Range("I1:J1").ColumnWidth = $wideFileColWidth

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $hoMdUrl, "", "", $hoMdDisplay)
$wsDeDe.Range("J2").Value = "6cf4087f-bfdd-4afa-b937-c2e3a4e86a9a.3dd1741d221c505698bb1d476ff36879f498027c.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-25 07:02:16"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $hoMdUrl, "", "", $hoMdDisplay)
$wsDeDe.Range("J3").Value = "6cf4087f-bfdd-4afa-b937-c2e3a4e86a9a.3dd1741d221c505698bb1d476ff36879f498027c.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-25 07:02:16"
